$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 25; existing rows 25-40 shift down to 27-42.
$ws.Rows("25:26").Insert()

# New row 25: Tuna "Especial", date 2023-03-16 (serial 45001)
$ws.Cells.Item(25,1).Value = 8
$ws.Cells.Item(25,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25,3).Value = "Coquimbo"
$ws.Cells.Item(25,4).Value = 45001
$ws.Cells.Item(25,5).Value = 4
$ws.Cells.Item(25,6).Value = "Fruta"
$ws.Cells.Item(25,7).Value = 100107
$ws.Cells.Item(25,8).Value = "Otros"
$ws.Cells.Item(25,9).Value = 100107011
$ws.Cells.Item(25,10).Value = "Tuna"
$ws.Cells.Item(25,11).Value = "Sin especificar"
$ws.Cells.Item(25,12).Value = "Especial"
$ws.Cells.Item(25,13).Value = 400
$ws.Cells.Item(25,14).Value = 12000
$ws.Cells.Item(25,15).Value = 13000
$ws.Cells.Item(25,16).Value = 12500
$ws.Cells.Item(25,17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(25,18).Value = "Provincia de Limarí"
$ws.Cells.Item(25,19).Value = 694
$ws.Cells.Item(25,20).Value = 18

# New row 26: Tuna "Primera", date 2023-03-16 (serial 45001)
$ws.Cells.Item(26,1).Value = 8
$ws.Cells.Item(26,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26,3).Value = "Coquimbo"
$ws.Cells.Item(26,4).Value = 45001
$ws.Cells.Item(26,5).Value = 4
$ws.Cells.Item(26,6).Value = "Fruta"
$ws.Cells.Item(26,7).Value = 100107
$ws.Cells.Item(26,8).Value = "Otros"
$ws.Cells.Item(26,9).Value = 100107011
$ws.Cells.Item(26,10).Value = "Tuna"
$ws.Cells.Item(26,11).Value = "Sin especificar"
$ws.Cells.Item(26,12).Value = "Primera"
$ws.Cells.Item(26,13).Value = 300
$ws.Cells.Item(26,14).Value = 10000
$ws.Cells.Item(26,15).Value = 11000
$ws.Cells.Item(26,16).Value = 10500
$ws.Cells.Item(26,17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(26,18).Value = "Provincia de Limarí"
$ws.Cells.Item(26,19).Value = 583
$ws.Cells.Item(26,20).Value = 18

Write-Host "Inserted rows 25-26; new dimension: $($ws.UsedRange.Address())"
